$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update footer timestamp text (row 1)
$ws.Range("A1").Value = "Datos actualizados a 5 de Mayo de 2020 a las 16:03"

# Row 4
$ws.Range("B4").Value = 1214023
$ws.Range("C4").Value = 1188
$ws.Range("D4").Value = 188069
$ws.Range("E4").Value = 955980
$ws.Range("G4").Value = 53
$ws.Range("H4").Value = 69974

# Row 30
$ws.Range("B30").Value = 18350
$ws.Range("C30").Value = 861
$ws.Range("D30").Value = 3771
$ws.Range("E30").Value = 14472
$ws.Range("G30").Value = 4
$ws.Range("H30").Value = 107

# Row 37
$ws.Range("E37").Value = 7542
$ws.Range("G37").Value = 23
$ws.Range("H37").Value = 841

# Row 44
$ws.Range("B44").Value = 9677
$ws.Range("C44").Value = 120
$ws.Range("D44").Value = 1723
$ws.Range("E44").Value = 7754
$ws.Range("F44").Value = 51
$ws.Range("G44").Value = 3
$ws.Range("H44").Value = 200

# Row 72
$ws.Range("D72").Value = 1501
$ws.Range("E72").Value = 693

# Row 75
$ws.Range("B75").Value = 2060
$ws.Range("C75").Value = 76
$ws.Range("D75").Value = 1508
$ws.Range("E75").Value = 526

# Row 99
$ws.Range("A99").Value = "Sri Lanka"
$ws.Range("B99").Value = 760
$ws.Range("C99").Value = 9
$ws.Range("D99").Value = 197
$ws.Range("E99").Value = 554
$ws.Range("F99").Value = 1
$ws.Range("G99").Value = 1
$ws.Range("H99").Value = 9

# Row 100
$ws.Range("A100").Value = "Somalia"
$ws.Range("B100").Value = 756
$ws.Range("C100").Value = 0
$ws.Range("D100").Value = 61
$ws.Range("E100").Value = 660
$ws.Range("F100").Value = 2
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 35

# Row 104
$ws.Range("D104").Value = 206
$ws.Range("E104").Value = 510

# Row 107
$ws.Range("A107").Value = "Consejo Danes para los Refugiados"
$ws.Range("B107").Value = 705
$ws.Range("C107").Value = 23
$ws.Range("D107").Value = 90
$ws.Range("E107").Value = 581
$ws.Range("F107").Value = 0
$ws.Range("H107").Value = 34

# Row 108
$ws.Range("A108").Value = "Mayotte"
$ws.Range("B108").Value = 686
$ws.Range("D108").Value = 352
$ws.Range("E108").Value = 328
$ws.Range("F108").Value = 6
$ws.Range("H108").Value = 6

# Row 111
$ws.Range("A111").Value = "Mali"
$ws.Range("B111").Value = 612
$ws.Range("C111").Value = 32
$ws.Range("D111").Value = 228
$ws.Range("E111").Value = 352
$ws.Range("F111").Value = 0
$ws.Range("G111").Value = 3
$ws.Range("H111").Value = 32

# Row 112
$ws.Range("A112").Value = "Georgia"
$ws.Range("B112").Value = 604
$ws.Range("C112").Value = 11
$ws.Range("D112").Value = 240
$ws.Range("E112").Value = 355
$ws.Range("F112").Value = 6
$ws.Range("H112").Value = 9

# Row 113
$ws.Range("A113").Value = "San Marino"
$ws.Range("B113").Value = 589
$ws.Range("C113").Value = 7
$ws.Range("D113").Value = 92
$ws.Range("E113").Value = 456
$ws.Range("F113").Value = 5
$ws.Range("H113").Value = 41

# Row 114
$ws.Range("A114").Value = "El Salvador"
$ws.Range("B114").Value = 587
$ws.Range("C114").Value = 32
$ws.Range("D114").Value = 201
$ws.Range("E114").Value = 373
$ws.Range("F114").Value = 3
$ws.Range("H114").Value = 13

# Row 126
$ws.Range("B126").Value = 366
$ws.Range("C126").Value = 4
$ws.Range("D126").Value = 127
$ws.Range("E126").Value = 237

# Row 128
$ws.Range("D128").Value = 319
$ws.Range("E128").Value = 3

# Row 132
$ws.Range("A132").Value = "Tayikistan"
$ws.Range("B132").Value = 293
$ws.Range("C132").Value = 63
$ws.Range("D132").Value = 0
$ws.Range("E132").Value = 288
$ws.Range("F132").Value = 0
$ws.Range("G132").Value = 2
$ws.Range("H132").Value = 5

# Row 133
$ws.Range("A133").Value = "Vietnam"
$ws.Range("B133").Value = 271
$ws.Range("D133").Value = 232
$ws.Range("E133").Value = 39
$ws.Range("F133").Value = 8

# Row 134
$ws.Range("A134").Value = "Ruanda"
$ws.Range("B134").Value = 261
$ws.Range("D134").Value = 128
$ws.Range("E134").Value = 133
$ws.Range("H134").Value = 0

# Row 135
$ws.Range("A135").Value = "Congo"
$ws.Range("B135").Value = 236
$ws.Range("D135").Value = 26
$ws.Range("E135").Value = 200
$ws.Range("H135").Value = 10

# Row 137
$ws.Range("A137").Value = "Cabo Verde"
$ws.Range("B137").Value = 186
$ws.Range("C137").Value = 11
$ws.Range("D137").Value = 37
$ws.Range("E137").Value = 147
$ws.Range("F137").Value = 0
$ws.Range("H137").Value = 2

# Row 138
$ws.Range("A138").Value = "Martinica"
$ws.Range("B138").Value = 181
$ws.Range("D138").Value = 83
$ws.Range("E138").Value = 84
$ws.Range("F138").Value = 5
$ws.Range("H138").Value = 14

# Row 139
$ws.Range("A139").Value = "Sierra Leona"
$ws.Range("B139").Value = 178
$ws.Range("E139").Value = 132
$ws.Range("H139").Value = 9

# Row 199
$ws.Range("A199").Value = "San Cristobal y Nieves"
$ws.Range("D199").Value = 8
$ws.Range("H199").Value = 0

# Row 200
$ws.Range("A200").Value = "Burundi"
$ws.Range("D200").Value = 7
$ws.Range("H200").Value = 1
